# This script reshapes "Sheet1" from a wide 10-column layout (with a
# duplicated / mis-cased date format) into the compact 5-column
# "Classement des éliminés" layout, matching the target commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Remember the existing bold/bordered/centered header style (currently
#    used on row 1, e.g. B1) by copying it, and immediately stamp it onto
#    the future header row (row 2) so the existing cellXfs entry gets
#    reused instead of a new font/xf combination being minted.
$ws.Range("B1").Copy()
$ws.Range("A2:E2").PasteSpecial(-4122)   # xlPasteFormats

# 2) Clear out everything else that needs to disappear: the old header
#    row (row 1), the extra columns F:J on rows 2-3, and the old row 3
#    content (we'll refill A3:D3 right after). Row A2:E2 is left alone so
#    the style stamped above survives.
$ws.Range("A1:J1").Clear()
$ws.Range("F2:J3").Clear()
$ws.Range("A3:E3").Clear()

# 3) Row 1: single title cell, no special style.
$ws.Range("A1").Value = "Classement des éliminés"

# 4) Row 2: the 5 header labels (style already applied above).
$ws.Range("A2").Value = "Classement"
$ws.Range("B2").Value = "Joueur"
$ws.Range("C2").Value = "Heure"
$ws.Range("D2").Value = "Killer"
$ws.Range("E2").Value = "Points"

# 5) Row 3: the actual data row.
$ws.Range("A3").Value = 0
$ws.Range("B3").Value = "Côme"
$ws.Range("C3").Value = 45542.70299954861
$ws.Range("C3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D3").Value = "Eric"

Write-Host "sheet rebuilt"
